$d = $word.ActiveDocument

# --- Hunk 1: merge "A" + " propriedade " into "A propriedade " ---
$d.Content.Find.Execute("A propriedade ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A propriedade ", 2)

# --- Hunk 2: merge " " + "B e transmite-a no sentido oposto. Como" ---
$d.Content.Find.Execute(" B e transmite-a no sentido oposto. Como", $true, $false, $false, $false, $false,
                         $true, 1, $false, " B e transmite-a no sentido oposto. Como", 2)

# --- Hunk 3: merge " é usar cookies com a " + "flag" + " " into " é usar cookies com a flag " ---
$d.Content.Find.Execute(" é usar cookies com a flag ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " é usar cookies com a flag ", 2)

# --- Hunk 4: append new bold sentence at end of exercise 3 paragraph ---
$target = "por código malicioso JavaScript de outros sites."
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*$target*") {
        $r = $p.Range
        $r.Collapse(0)
        $r.InsertAfter(" ")
        $r.Collapse(0)
        $r.InsertAfter("Interpretaste mal a pergunta.")
        $r.Bold = 1
        break
    }
}
